# Edit: update the Hortaliza / Comercializadora del Agro de Limarí - Pepino
# ensalada dataset for the "Fruta / hortaliza, semanal" weekly refresh.
#
# The underlying rows (2..90) get reshuffled/updated with newer market-report
# values (dates, volumes, prices, unit of sale, $/Kg, etc.), and two brand
# new observations are appended as rows 91 and 92. This mirrors the upstream
# commit's unified diff, row by row, column by column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference number format for the "Fecha" (date) column, taken from an
# existing data row, so newly-created rows 91/92 keep the same date format.
$dateFormat = $ws.Cells.Item(2, 4).NumberFormat

$data = @(
    @{Row=2; D=44231; I="Primera"; J=700; K=9000; L=9500; M=9250; N="$/caja 70 unidades"; P=132; Q=70}
    @{Row=3; D=44294; I="Primera"; J=700; K=13500; L=14000; M=13750; N="$/caja 70 unidades"; P=196; Q=70}
    @{Row=4; D=44294; I="Segunda"; J=240; K=10500; L=11000; M=10750; N="$/caja 100 unidades"; P=108; Q=100}
    @{Row=5; D=44223; I="Primera"; J=700; K=9500; L=10000; M=9750; N="$/caja 70 unidades"; P=139; Q=70}
    @{Row=6; D=44223; I="Segunda"; J=400; K=7500; L=8000; M=7750; N="$/caja 100 unidades"; P=78; Q=100}
    @{Row=7; D=44371; I="Primera"; J=1000; K=12500; L=13000; M=12750; N="$/caja 60 unidades"; P=212; Q=60}
    @{Row=8; D=44371; I="Segunda"; J=400; K=10500; L=11000; M=10750; N="$/caja 100 unidades"; P=108; Q=100}
    @{Row=9; D=44503; I="Primera"; J=1000; K=6500; L=7000; M=6750; N="$/caja 70 unidades"; P=96; Q=70}
    @{Row=10; D=44503; I="Segunda"; J=1100; K=4500; L=5000; M=4750; N="$/caja 100 unidades"; P=48; Q=100}
    @{Row=11; D=44258; I="Primera"; J=1000; K=11000; L=12000; M=11500; N="$/caja 70 unidades"; P=164; Q=70}
    @{Row=12; D=44258; I="Segunda"; J=500; K=8000; L=9000; M=8500; N="$/caja 100 unidades"; P=85; Q=100}
    @{Row=13; D=44160; I="Primera"; J=1700; K=6500; L=7000; M=6750; N="$/caja 70 unidades"; P=96; Q=70}
    @{Row=14; D=44385; I="Primera"; J=500; K=13000; L=14000; M=13500; N="$/caja 60 unidades"; P=225; Q=60}
    @{Row=15; D=44385; I="Segunda"; J=360; K=11000; L=12000; M=11500; N="$/caja 100 unidades"; P=115; Q=100}
    @{Row=16; D=44370; I="Primera"; J=800; K=12500; L=13000; M=12750; N="$/caja 60 unidades"; P=212; Q=60}
    @{Row=17; D=44370; I="Segunda"; J=500; K=10500; L=11000; M=10750; N="$/caja 100 unidades"; P=108; Q=100}
    @{Row=18; D=44280; I="Primera"; J=600; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=19; D=44280; I="Segunda"; J=240; K=8500; L=9000; M=8750; N="$/caja 100 unidades"; P=88; Q=100}
    @{Row=20; D=44377; I="Primera"; J=1000; K=13000; L=14000; M=13500; N="$/caja 60 unidades"; P=225; Q=60}
    @{Row=21; D=44377; I="Segunda"; J=500; K=9000; L=10000; M=9500; N="$/caja 100 unidades"; P=95; Q=100}
    @{Row=22; D=44188; I="Primera"; J=2500; K=6500; L=7000; M=6750; N="$/caja 70 unidades"; P=96; Q=70}
    @{Row=23; D=44230; I="Primera"; J=1100; K=9000; L=10000; M=9500; N="$/caja 70 unidades"; P=136; Q=70}
    @{Row=24; D=44266; I="Primera"; J=600; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=25; D=44224; I="Primera"; J=1100; K=9500; L=10000; M=9750; N="$/caja 70 unidades"; P=139; Q=70}
    @{Row=26; D=44224; I="Segunda"; J=600; K=7500; L=8000; M=7750; N="$/caja 100 unidades"; P=78; Q=100}
    @{Row=27; D=44286; I="Primera"; J=1200; K=11000; L=12000; M=11500; N="$/caja 70 unidades"; P=164; Q=70}
    @{Row=28; D=44335; I="Primera"; J=700; K=10000; L=11000; M=10500; N="$/caja 60 unidades"; P=175; Q=60}
    @{Row=29; D=44335; I="Segunda"; J=500; K=8000; L=9000; M=8500; N="$/caja 100 unidades"; P=85; Q=100}
    @{Row=30; D=44392; I="Segunda"; J=300; K=14000; L=15000; M=14500; N="$/caja 100 unidades"; P=145; Q=100}
    @{Row=31; D=44497; I="Primera"; J=600; K=6500; L=7000; M=6750; N="$/caja 70 unidades"; P=96; Q=70}
    @{Row=32; D=44497; I="Segunda"; J=680; K=4500; L=5000; M=4721; N="$/caja 100 unidades"; P=47; Q=100}
    @{Row=33; D=44273; I="Primera"; J=500; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=34; D=44273; I="Segunda"; J=360; K=8000; L=9000; M=8500; N="$/caja 100 unidades"; P=85; Q=100}
    @{Row=35; D=44175; I="Primera"; J=2400; K=6500; L=7000; M=6750; N="$/caja 70 unidades"; P=96; Q=70}
    @{Row=36; D=44175; I="Segunda"; J=1700; K=4500; L=5000; M=4750; N="$/caja 100 unidades"; P=48; Q=100}
    @{Row=37; D=44168; I="Primera"; J=1700; K=6000; L=7000; M=6500; N="$/caja 60 unidades"; P=108; Q=60}
    @{Row=38; D=44203; I="Primera"; J=500; K=8500; L=9000; M=8750; N="$/caja 70 unidades"; P=125; Q=70}
    @{Row=39; D=44244; I="Primera"; J=1600; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=40; D=44244; I="Segunda"; J=500; K=7000; L=8000; M=7500; N="$/caja 100 unidades"; P=75; Q=100}
    @{Row=41; D=44202; I="Primera"; J=400; K=8500; L=9000; M=8750; N="$/caja 70 unidades"; P=125; Q=70}
    @{Row=42; D=44252; I="Primera"; J=1100; K=9000; L=10000; M=9500; N="$/caja 70 unidades"; P=136; Q=70}
    @{Row=43; D=44252; I="Segunda"; J=800; K=7500; L=8000; M=7750; N="$/caja 100 unidades"; P=78; Q=100}
    @{Row=44; D=44167; I="Primera"; J=1700; K=6000; L=7000; M=6500; N="$/caja 60 unidades"; P=108; Q=60}
    @{Row=45; D=44167; I="Segunda"; J=1000; K=4000; L=5000; M=4500; N="$/caja 100 unidades"; P=45; Q=100}
    @{Row=46; D=44293; I="Primera"; J=1000; K=13500; L=14000; M=13750; N="$/caja 70 unidades"; P=196; Q=70}
    @{Row=47; D=44293; I="Segunda"; J=240; K=10500; L=11000; M=10750; N="$/caja 100 unidades"; P=108; Q=100}
    @{Row=48; D=44308; I="Primera"; J=600; K=9500; L=10000; M=9750; N="$/caja 70 unidades"; P=139; Q=70}
    @{Row=49; D=44308; I="Segunda"; J=400; K=7500; L=8000; M=7750; N="$/caja 100 unidades"; P=78; Q=100}
    @{Row=50; D=44321; I="Primera"; J=500; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=51; D=44189; I="Primera"; J=700; K=7000; L=7500; M=7250; N="$/caja 70 unidades"; P=104; Q=70}
    @{Row=52; D=44265; I="Primera"; J=500; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=53; D=44300; I="Primera"; J=1000; K=12000; L=13000; M=12500; N="$/caja 70 unidades"; P=179; Q=70}
    @{Row=54; D=44300; I="Segunda"; J=400; K=9000; L=10000; M=9500; N="$/caja 100 unidades"; P=95; Q=100}
    @{Row=55; D=44209; I="Primera"; J=700; K=7500; L=8000; M=7750; N="$/caja 70 unidades"; P=111; Q=70}
    @{Row=56; D=44209; I="Segunda"; J=500; K=5500; L=6000; M=5750; N="$/caja 100 unidades"; P=58; Q=100}
    @{Row=57; D=44237; I="Primera"; J=600; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=58; D=44237; I="Segunda"; J=700; K=7000; L=8000; M=7500; N="$/caja 100 unidades"; P=75; Q=100}
    @{Row=59; D=44259; I="Primera"; J=1000; K=11000; L=12000; M=11500; N="$/caja 70 unidades"; P=164; Q=70}
    @{Row=60; D=44259; I="Segunda"; J=400; K=8000; L=9000; M=8500; N="$/caja 100 unidades"; P=85; Q=100}
    @{Row=61; D=44363; I="Primera"; J=1000; K=11000; L=12000; M=11500; N="$/caja 60 unidades"; P=192; Q=60}
    @{Row=62; D=44363; I="Segunda"; J=700; K=9000; L=10000; M=9500; N="$/caja 100 unidades"; P=95; Q=100}
    @{Row=63; D=44336; I="Primera"; J=500; K=10000; L=11000; M=10500; N="$/caja 60 unidades"; P=175; Q=60}
    @{Row=64; D=44336; I="Segunda"; J=400; K=7000; L=8000; M=7500; N="$/caja 100 unidades"; P=75; Q=100}
    @{Row=65; D=44195; I="Primera"; J=800; K=7500; L=8000; M=7750; N="$/caja 70 unidades"; P=111; Q=70}
    @{Row=66; D=44195; I="Segunda"; J=400; K=5500; L=6000; M=5750; N="$/caja 100 unidades"; P=58; Q=100}
    @{Row=67; D=44210; I="Primera"; J=500; K=7500; L=8000; M=7750; N="$/caja 70 unidades"; P=111; Q=70}
    @{Row=68; D=44210; I="Segunda"; J=400; K=5500; L=6000; M=5750; N="$/caja 100 unidades"; P=58; Q=100}
    @{Row=69; D=44301; I="Primera"; J=500; K=12000; L=13000; M=12500; N="$/caja 70 unidades"; P=179; Q=70}
    @{Row=70; D=44301; I="Segunda"; J=400; K=9000; L=10000; M=9500; N="$/caja 100 unidades"; P=95; Q=100}
    @{Row=71; D=44279; I="Primera"; J=700; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=72; D=44279; I="Segunda"; J=500; K=8500; L=9000; M=8750; N="$/caja 100 unidades"; P=88; Q=100}
    @{Row=73; D=44504; I="Primera"; J=800; K=6500; L=7000; M=6750; N="$/caja 70 unidades"; P=96; Q=70}
    @{Row=74; D=44504; I="Segunda"; J=700; K=4500; L=5000; M=4750; N="$/caja 100 unidades"; P=48; Q=100}
    @{Row=75; D=44384; I="Primera"; J=400; K=13000; L=14000; M=13500; N="$/caja 60 unidades"; P=225; Q=60}
    @{Row=76; D=44384; I="Segunda"; J=500; K=11000; L=12000; M=11500; N="$/caja 100 unidades"; P=115; Q=100}
    @{Row=77; D=44272; I="Primera"; J=400; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=78; D=44272; I="Segunda"; J=400; K=8000; L=9000; M=8500; N="$/caja 100 unidades"; P=85; Q=100}
    @{Row=79; D=44322; I="Primera"; J=400; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=80; D=44510; I="Primera"; J=400; K=5500; L=6000; M=5750; N="$/caja 70 unidades"; P=82; Q=70}
    @{Row=81; D=44161; I="Primera"; J=1100; K=6000; L=6500; M=6250; N="$/caja 70 unidades"; P=89; Q=70}
    @{Row=82; D=44517; I="Primera"; J=700; K=5500; L=6000; M=5750; N="$/caja 70 unidades"; P=82; Q=70}
    @{Row=83; D=44517; I="Segunda"; J=600; K=3500; L=4000; M=3750; N="$/caja 100 unidades"; P=38; Q=100}
    @{Row=84; D=44238; I="Primera"; J=700; K=10000; L=11000; M=10500; N="$/caja 70 unidades"; P=150; Q=70}
    @{Row=85; D=44238; I="Segunda"; J=600; K=7000; L=8000; M=7500; N="$/caja 100 unidades"; P=75; Q=100}
    @{Row=86; D=44391; I="Primera"; J=400; K=14000; L=15000; M=14500; N="$/caja 60 unidades"; P=242; Q=60}
    @{Row=87; D=44391; I="Segunda"; J=240; K=12000; L=13000; M=12500; N="$/caja 100 unidades"; P=125; Q=100}
    @{Row=88; D=44251; I="Primera"; J=1200; K=9000; L=10000; M=9500; N="$/caja 70 unidades"; P=136; Q=70}
    @{Row=89; D=44251; I="Segunda"; J=700; K=7000; L=8000; M=7500; N="$/caja 100 unidades"; P=75; Q=100}
    @{Row=90; D=44181; I="Primera"; J=600; K=6000; L=6500; M=6250; N="$/caja 70 unidades"; P=89; Q=70}
    @{Row=91; D=44307; I="Primera"; J=500; K=9000; L=10000; M=9500; N="$/caja 70 unidades"; P=136; Q=70}
    @{Row=92; D=44307; I="Segunda"; J=400; K=7000; L=8000; M=7500; N="$/caja 100 unidades"; P=75; Q=100}
)

foreach ($item in $data) {
    $r = $item.Row

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $item.D
    $dCell.NumberFormat = $dateFormat

    $ws.Cells.Item($r, 9).Value  = $item.I    # Calidad
    $ws.Cells.Item($r, 10).Value = $item.J    # Volumen
    $ws.Cells.Item($r, 11).Value = $item.K    # Precio minimo
    $ws.Cells.Item($r, 12).Value = $item.L    # Precio maximo
    $ws.Cells.Item($r, 13).Value = $item.M    # Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $item.N    # Unidad de comercializacion
    $ws.Cells.Item($r, 16).Value = $item.P    # Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $item.Q    # Kg o Unidades
}

# New rows 91/92 don't yet carry the constant columns shared by every
# observation in this sheet (Mercado ID, Mercado, Región, Codreg, Categoría
# ID, Categoría, Variedad, Origen, Clasificación) - copy them from row 90.
foreach ($r in 91, 92) {
    $ws.Cells.Item($r, 1).Value  = $ws.Cells.Item(90, 1).Value2   # Mercado ID
    $ws.Cells.Item($r, 2).Value  = $ws.Cells.Item(90, 2).Value2   # Mercado
    $ws.Cells.Item($r, 3).Value  = $ws.Cells.Item(90, 3).Value2   # Región
    $ws.Cells.Item($r, 5).Value  = $ws.Cells.Item(90, 5).Value2   # Codreg
    $ws.Cells.Item($r, 6).Value  = $ws.Cells.Item(90, 6).Value2   # Categoría ID
    $ws.Cells.Item($r, 7).Value  = $ws.Cells.Item(90, 7).Value2   # Categoría
    $ws.Cells.Item($r, 8).Value  = $ws.Cells.Item(90, 8).Value2   # Variedad
    $ws.Cells.Item($r, 15).Value = $ws.Cells.Item(90, 15).Value2  # Origen
    $ws.Cells.Item($r, 18).Value = $ws.Cells.Item(90, 18).Value2  # Clasificación
}
